# normalize club names in ICL data and run south coast r1 file
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full refreshed ladder data (rows 2-27), reflecting the south coast round 1
# results merge: new clubs added (Illawarra, Kiama, Eurocoast, Jervis Bay)
# and the ladder re-sorted by Total Points.
$data = @(
    @("Warringah Triathlon Club", 30, 109, 139, 139, 253),
    @("Balmoral Triathlon Club", 30, 93, 123, 123, 152),
    @("Illawarra Triathlon Club", 30, 86, 116, 116, 77),
    @("Panthers Triathlon Club", 30, 79, 109, 109, 135),
    @("Moore Performance Triathlon Club", 45, 54, 99, 99, 44),
    @("STG Triathlon Club", 45, 47, 92, 92, 49),
    @("Kiama Triathlon Club", 30, 59, 89, 89, 67),
    @("Hills Triathlon Club", 30, 54, 84, 84, 178),
    @("Cronulla Triathlon Club", 15, 65, 80, 80, 207),
    @("Balance Triathlon Club", 30, 43, 73, 73, 127),
    @("Eurocoast Triathlon Club", 45, 26, 71, 71, 25),
    @("Coogee Triathlon Club", 30, 40, 70, 70, 119),
    @("Pulse Performance", 30, 26, 56, 56, 31),
    @("Concord Triathlon Club", 30, 15, 45, 45, 53),
    @("Engadine Triathlon Club", 45, 0, 45, 45, 8),
    @("Manly Vipers Triathlon Club", 30, 9, 39, 39, 33),
    @("BRAT Triathlon Club", 30, 9, 39, 39, 161),
    @("Macarthur Triathlon Club", 15, 19, 34, 34, 69),
    @("Jervis Bay Triathlon Club", 15, 16, 31, 31, 47),
    @("FilOz Triathlon Club", 30, 0, 30, 30, 10),
    @("Sydney South West Triathlon Club", 30, 0, 30, 30, 2),
    @("Brighton Baths Athletic Club", 15, 9, 24, 24, 42),
    @("Northern Suburbs Triathlon Club", 15, 9, 24, 24, 67),
    @("Australian Chinese Triathlon Club", 0, 0, 0, 0, 3),
    @("Hunters Hills Triathlon Club", 0, 0, 0, 0, 6),
    @("TriMob", 0, 0, 0, 0, 43)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $ws.Cells.Item($row, 5).Value = $entry[4]
    $ws.Cells.Item($row, 6).Value = $entry[5]
    $row = $row + 1
}
